# Remove the "ram"/"people" employee record (was row 3 in the sheet, i.e.
# worksheet row 3) and shift the following rows up, matching the commit
# that removed this row's SQL-server-driven data manipulation code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:3").Delete()
